$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projeto")

# Update language column (C2:C24): remove surrounding quotes from the text values
$ws.Range("C2:C22").Value = "pt"
$ws.Range("C23").Value = "es"
$ws.Range("C24").Value = "en"

# Update the active selection to B1 (also resets the frozen pane's
# top-left visible cell back to A2)
$ws.Activate()
$ws.Range("B1").Select()

# Adjust column B width (no longer best-fit / auto-fit)
$ws.Columns("B").ColumnWidth = 24.5

$wb.Save()
